# Update the "type" (column B) and "value" (column C) columns for all
# Flag Football rows (rows 2-28) on the active sheet.
#
# Previously:
#   B = "club-sports"
#   C = "Flag Football-Girls" | "Flag Football-Coed" | "Flag Football-Boys"
#
# Now:
#   B = "sports_club_girls" | "sports_club_coed" | "sports_club_boys"
#       (derived from what C used to contain)
#   C = "Flag Football" (suffix removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $currentValue = $ws.Cells.Item($r, 3).Value2

    if ($currentValue -eq "Flag Football-Girls") {
        $newType = "sports_club_girls"
    } elseif ($currentValue -eq "Flag Football-Coed") {
        $newType = "sports_club_coed"
    } elseif ($currentValue -eq "Flag Football-Boys") {
        $newType = "sports_club_boys"
    } else {
        $newType = $null
    }

    if ($newType -ne $null) {
        $ws.Cells.Item($r, 2).Value = $newType
        $ws.Cells.Item($r, 3).Value = "Flag Football"
    }
}
